# Update "DateBase/orders/Dang Nguyen 195_2026-2-9.xlsx"
#  - Orders sheet: append order-line rows 67-73 (extends used range to A1:L73)
#  - Summary sheet: extend the packed G2 tracking string

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Orders" (sheet1) — append rows 67-73
# ---------------------------------------------------------------------------
$orders = $wb.Worksheets.Item("Orders")

# Helper cells store plain numbers as TEXT (matches the rest of the sheet,
# e.g. A2/F2 etc. are all text-typed), so force the "Text" number format on
# the numeric-looking cells before writing the values (only the cells that
# actually receive a value - leave the rest of column A alone since most of
# the new rows have no PackageID).
$orders.Range("A67:A68").NumberFormat = "@"
$orders.Range("A71").NumberFormat = "@"
$orders.Range("F67:F73").NumberFormat = "@"

$newline = [char]10

$orders.Cells.Item(67, 1).Value = "1"
$orders.Cells.Item(67, 3).Value = "775_海芋黑_Calla Lily_undefined_1bunch"
$orders.Cells.Item(67, 6).Value = "8"

$orders.Cells.Item(68, 1).Value = "2"
$orders.Cells.Item(68, 3).Value = "653_大丽花 黑_undefined_undefined_5stems"
$orders.Cells.Item(68, 6).Value = "16"

$orders.Cells.Item(69, 3).Value = "592_进口春兰叶_undefined_undefined_1bunch"
$orders.Cells.Item(69, 6).Value = "5"

$orders.Cells.Item(70, 3).Value = "344_钢草_steal grass_Xanthorrhoea preissii Endl._1bunch"
$orders.Cells.Item(70, 6).Value = "5"

$orders.Cells.Item(71, 1).Value = "3"
$orders.Cells.Item(71, 3).Value = [string]::Concat("542_吊米 红_hanging amaranthus", $newline, "red_undefined_1bunch")
$orders.Cells.Item(71, 6).Value = "5"

$orders.Cells.Item(72, 3).Value = "484_天鹅绒_Star of Bethlehem_undefined_1bunch"
$orders.Cells.Item(72, 6).Value = "5"

$orders.Cells.Item(73, 3).Value = "3_波浪白洋桔梗_Wavy White Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
$orders.Cells.Item(73, 6).Value = "4"

# ---------------------------------------------------------------------------
# Sheet "Summary" (sheet2) — extend the packed tracking string in G2
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$g2 = $summary.Cells.Item(2, 7)
$g2old = $g2.Text
$g2.NumberFormat = "@"
$g2.Value = [string]::Concat($g2old, "81655554")
